$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D hold text values that look numeric (e.g. "0.9990",
# "28.813.47"). Force them to Text format before assigning so Excel
# does not silently reinterpret them as numbers (which would drop
# significant trailing zeros / thousand-separator dots). The style is
# reset back to Normal afterwards so no stray formatting is introduced.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.813.47'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +7.43%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.814.45'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +5.12%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9990'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.12%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '250.96'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +4.03%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9990'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.09%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4985'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.44%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2787'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +7.75%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06384'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.95%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.814.76'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +5.12%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '16.80'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +5.38%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07164'
$ws.Range("D12").Style = "Normal"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.6499'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +7.25%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.713'
$ws.Range("D14").Style = "Normal"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '82.26'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +6.76%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '28.794.46'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +8.24%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.9991'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.05%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000007395'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +3.39%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.9987'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.09%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.30'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +7.60%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.048.43'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +4.91%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.619'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +4.55%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.900'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +3.96%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.353'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +5.48%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '144.04'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +4.67%  '

$ws.Range("B26").Value = 'BitcoinCash'
$ws.Range("C26").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '122.38'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +15.08%  '

$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '16.17'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +6.11%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.897'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +7.09%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.397'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.34%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.189'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +6.57%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08382'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +5.11%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.861'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +5.00%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04972'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +10.37%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.091'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +8.05%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.6818'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +9.34%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.683'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +3.20%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.743'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +12.99%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.9689'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +4.41%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.195'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +8.11%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01591'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +6.42%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.020'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +6.41%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9993'
$ws.Range("D42").Style = "Normal"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '101.72'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.51%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.4120'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +7.55%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '7.252'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +5.93%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.1227'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +6.03%  '

$ws.Range("E47").Value = '  +2.00%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.181'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +4.14%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '31.72'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +5.32%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.3654'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +8.45%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.312'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +6.54%  '

